$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.156176328659058
$ws.Range("B1").Value = 5.946654796600342
$ws.Range("C1").Value = 1.233901500701904
$ws.Range("D1").Value = 0.2465873062610626
$ws.Range("E1").Value = 0.3340835273265839
